$d = $word.ActiveDocument

# Merge the three runs "<id>", "p052r_1", "</id>" into a single run whose
# formatting matches the first run ("<id>") by replacing the combined text
# in-place. Word's Find/Replace collapses the matched range into one run
# that carries the formatting of the first matched run.
$null = $d.Content.Find.Execute("<id>p052r_1</id>", $true, $false, $false, `
    $false, $false, $true, 1, $false, "<id>p052r_1</id>", 2)
